# Scheduled-runner refresh: re-pull current market-board prices and
# recompute the dependent Leve profit columns (H:N) for each job sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 3666.3076
$ws.Range("I31").Value = 537.2
$ws.Range("J31").Value = 5622
$ws.Range("K31").Value = 1611.6
$ws.Range("L31").Value = 16866
$ws.Range("M31").Value = -1381.6
$ws.Range("N31").Value = -17326

$ws.Range("H129").Value = 1163.3572
$ws.Range("I129").Value = 360.33334
$ws.Range("J129").Value = 1317.1277
$ws.Range("K129").Value = 1081.00002
$ws.Range("L129").Value = 3951.3831
$ws.Range("M129").Value = 3918.99998
$ws.Range("N129").Value = -13951.3831

$ws.Range("H137").Value = 28121.027
$ws.Range("I137").Value = 35364.758
$ws.Range("J137").Value = 1862.5
$ws.Range("K137").Value = 106094.274
$ws.Range("L137").Value = 5587.5
$ws.Range("M137").Value = -103544.274
$ws.Range("N137").Value = -10687.5

$ws.Range("H141").Value = 1682.878
$ws.Range("I141").Value = 1594.5
$ws.Range("J141").Value = 1704.303
$ws.Range("K141").Value = 4783.5
$ws.Range("L141").Value = 5112.909000000001
$ws.Range("M141").Value = 396.5
$ws.Range("N141").Value = -15472.909

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 232
$ws.Range("I5").Value = 232
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 232
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -120
$ws.Range("N5").ClearContents()

$ws.Range("H12").Value = 1900
$ws.Range("I12").Value = 1900
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 1900
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -1727

$ws.Range("H19").Value = 1000
$ws.Range("I19").Value = 1000
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 1000
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -771

$ws.Range("H32").Value = 27229.414
$ws.Range("I32").Value = 4966.851
$ws.Range("J32").Value = 122351.27
$ws.Range("K32").Value = 4966.851
$ws.Range("L32").Value = 122351.27
$ws.Range("M32").Value = -4679.851
$ws.Range("N32").Value = -122925.27

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 232
$ws.Range("I4").Value = 232
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 232
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -117
$ws.Range("N4").ClearContents()

$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

$ws.Range("H88").Value = 27743
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 27743
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 27743
$ws.Range("N88").Value = -28555

$ws.Range("H91").Value = 27743
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 27743
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 27743
$ws.Range("N91").Value = -30551

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()

$ws.Range("H22").Value = 415.625
$ws.Range("I22").Value = 446.05264
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 446.05264
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = -96.05264
$ws.Range("N22").Value = -1000

$ws.Range("H31").Value = 20001906
$ws.Range("I31").Value = 47620412
$ws.Range("J31").Value = 2297.8965
$ws.Range("K31").Value = 47620412
$ws.Range("L31").Value = 2297.8965
$ws.Range("M31").Value = -47620117
$ws.Range("N31").Value = -2887.8965

$ws.Range("H34").Value = 20001906
$ws.Range("I34").Value = 47620412
$ws.Range("J34").Value = 2297.8965
$ws.Range("K34").Value = 47620412
$ws.Range("L34").Value = 2297.8965
$ws.Range("M34").Value = -47620210
$ws.Range("N34").Value = -2701.8965

$ws.Range("H132").Value = 2158.9546
$ws.Range("I132").Value = 1247.8462
$ws.Range("J132").Value = 3475
$ws.Range("K132").Value = 3743.5386
$ws.Range("L132").Value = 10425
$ws.Range("M132").Value = -1213.5386
$ws.Range("N132").Value = -15485

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 721
$ws.Range("I12").Value = 270
$ws.Range("J12").Value = 811.2
$ws.Range("K12").Value = 810
$ws.Range("L12").Value = 2433.6
$ws.Range("M12").Value = -637
$ws.Range("N12").Value = -2779.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 1450.5
$ws.Range("I41").Value = 1001
$ws.Range("J41").Value = 1900
$ws.Range("K41").Value = 1001
$ws.Range("L41").Value = 1900
$ws.Range("M41").Value = -646

$ws.Range("H132").Value = 2227.3713
$ws.Range("I132").Value = 1644.8889
$ws.Range("J132").Value = 4193.25
$ws.Range("K132").Value = 4934.6667
$ws.Range("L132").Value = 12579.75
$ws.Range("M132").Value = -2404.6667
$ws.Range("N132").Value = -17639.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 15000
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 15000
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 15000
$ws.Range("N18").Value = -15344

$ws.Range("H32").Value = 2480.6667
$ws.Range("I32").Value = 2480.6667
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 2480.6667
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -2163.6667
$ws.Range("N32").ClearContents()

$ws.Range("H40").Value = 1967.7273
$ws.Range("I40").Value = 1721.4286
$ws.Range("J40").Value = 2398.75
$ws.Range("K40").Value = 1721.4286
$ws.Range("L40").Value = 2398.75
$ws.Range("M40").Value = -1585.4286
$ws.Range("N40").Value = -2670.75

$ws.Range("H55").Value = 917.4167
$ws.Range("I55").Value = 448
$ws.Range("J55").Value = 1152.125
$ws.Range("K55").Value = 448
$ws.Range("L55").Value = 1152.125
$ws.Range("M55").Value = -275
$ws.Range("N55").Value = -1498.125

$ws.Range("H81").Value = 31000
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 31000
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 31000
$ws.Range("N81").Value = -32996

$ws.Range("H84").Value = 31000
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 31000
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 93000
$ws.Range("N84").Value = -102984

$ws.Range("H122").Value = 3576.7856
$ws.Range("I122").Value = 2928.5715
$ws.Range("J122").Value = 4225
$ws.Range("K122").Value = 8785.7145
$ws.Range("L122").Value = 12675
$ws.Range("M122").Value = -6335.7145
$ws.Range("N122").Value = -17575

$ws.Range("H130").Value = 40819.168
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 40819.168
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 40819.168
$ws.Range("N130").Value = -50859.168

$ws.Range("H132").Value = 2738.75
$ws.Range("I132").Value = 1982.4
$ws.Range("J132").Value = 3999.3333
$ws.Range("K132").Value = 5947.200000000001
$ws.Range("L132").Value = 11997.9999
$ws.Range("M132").Value = -3417.200000000001
$ws.Range("N132").Value = -17057.9999

Write-Host "Sheets refreshed: ALC, ARM, BSM, CRP, CUL, GSM, LTW"
